# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Only column G ("K") values change for rows 2-25 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed "K" (column G) values per row, replacing the old Strike# figures.
$newK = @{
    2  = 3
    3  = 2
    4  = 0
    5  = 1
    6  = 1
    7  = 3
    8  = 2
    9  = 3
    10 = 1
    11 = 5
    12 = 1
    13 = 3
    14 = 1
    15 = 7
    16 = 9
    17 = 7
    18 = 6
    19 = 5
    20 = 4
    21 = 1
    22 = 5
    23 = 5
    24 = 3
    25 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
